$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The diff removes the old row 3 (id 252456) entirely, shifting later rows up.
$ws.Rows(3).Delete()

# After the shift, the old row 5 (id 251790) is now row 4; remove it too so the
# old row 6 (id 251231) becomes the new row 4.
$ws.Rows(4).Delete()

# Finally, everything from (the now) row 5 through row 16 is no longer present
# in the target sheet, so remove those trailing rows.
$ws.Range("5:16").EntireRow.Delete()
